# Add a new "student_id / staff_id" column (I) to the Grade4 roster sheet.
# Each data row (2-18) gets an ID value in column I, formatted the same way
# as the rest of that row (copy format from column B of the same row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @{
    2  = "16-0175"
    3  = "17-0057"
    4  = "18-0058"
    5  = "16-0057"
    6  = "17-0164"
    7  = "17-0036"
    8  = "17-0088"
    9  = "16-0084"
    10 = "16-0091"
    11 = "17-0033"
    12 = "18-0213"
    13 = "16-0126"
    14 = "18-0006"
    15 = "18-0143"
    16 = "16-0128"
    17 = "16-0143"
    18 = "16-0173"
}

foreach ($row in 2..18) {
    # Match the formatting already used for the rest of the row (column B
    # carries the row's "data" style).
    $ws.Range("B$row").Copy()
    $ws.Range("I$row").PasteSpecial(-4122)

    $ws.Range("I$row").Value = $ids[$row]
}

# Reflect the new column in the view: select it and scroll back to A1, same
# as Excel would leave things after adding/filling the column.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 90
$ws.Range("I2:I18").Select()

# Slightly narrower columns (consistent with the new zoom level).
$ws.Columns.Item(1).ColumnWidth = 12.833333333333334
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334
$ws.Columns.Item(3).ColumnWidth = 23.5
$ws.Columns.Item(4).ColumnWidth = 12.833333333333334
$ws.Columns.Item(5).ColumnWidth = 32.833333333333336
$ws.Columns.Item(6).ColumnWidth = 33.5
$ws.Columns.Item(7).ColumnWidth = 29.5
$ws.Columns.Item(8).ColumnWidth = 7.5
